# Update gh-pages to output generated at 456a3b4
# Refresh "想去人数" (column F, interest-count) figures scraped from bilibili,
# and one stale cover-image URL (I41 on 展览), across all four sheets.

$wb = $excel.ActiveWorkbook

# --- 展览 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 162
$ws.Range("F3").Value = 936
$ws.Range("F4").Value = 1095
$ws.Range("F5").Value = 1555
$ws.Range("F7").Value = 698
$ws.Range("F8").Value = 12637
$ws.Range("F9").Value = 2220
$ws.Range("F13").Value = 37250
$ws.Range("F14").Value = 1258
$ws.Range("F15").Value = 249
$ws.Range("F16").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("F21").Value = 800
$ws.Range("F22").Value = 4579
$ws.Range("F23").Value = 4580
$ws.Range("F24").Value = 1175
$ws.Range("F25").Value = 896
$ws.Range("F28").Value = 22
$ws.Range("F29").Value = 7
$ws.Range("F30").Value = 1122
$ws.Range("F31").Value = 65
$ws.Range("F32").Value = 128
$ws.Range("F33").Value = 290
$ws.Range("F36").Value = 42
$ws.Range("F38").Value = 4532
$ws.Range("F39").Value = 22
$ws.Range("F40").Value = 4635
$ws.Range("F41").Value = 5594
$ws.Range("F43").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("F46").Value = 377
$ws.Range("F47").Value = 88
$ws.Range("F48").Value = 55
$ws.Range("F49").Value = 4128
$ws.Range("I41").Value = "//i2.hdslb.com/bfs/openplatform/202409/6ZZ2zi6T1727663088882.jpeg"

# --- 演出 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F5").Value = 108
$ws.Range("F8").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("F16").Value = 6
$ws.Range("F19").Value = 0

# --- 本地生活 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value = 494
$ws.Range("F4").Value = 113
$ws.Range("F5").Value = 0

# --- 全部类型 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 494
$ws.Range("F3").Value = 162
$ws.Range("F4").Value = 936
$ws.Range("F5").Value = 1095
$ws.Range("F7").Value = 698
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = 12637
$ws.Range("F10").Value = 2220
$ws.Range("F13").Value = 1258
$ws.Range("F14").Value = 285
$ws.Range("F16").Value = 0
$ws.Range("F17").Value = 326
$ws.Range("F19").Value = 800
$ws.Range("F21").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("F26").Value = 896
$ws.Range("F27").Value = 0
$ws.Range("F28").Value = 7
$ws.Range("F29").Value = 1122
$ws.Range("F30").Value = 0
$ws.Range("F31").Value = 65
$ws.Range("F32").Value = 128
$ws.Range("F34").Value = 290
$ws.Range("F37").Value = 0
$ws.Range("F38").Value = 4635
$ws.Range("F39").Value = 0
$ws.Range("F41").Value = 98
$ws.Range("F43").Value = 377
$ws.Range("F44").Value = 6
$ws.Range("F46").Value = 88
$ws.Range("F47").Value = 4128
